$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2789.4736
$ws.Range("I112").Value = 816.6667
$ws.Range("J112").Value = 3159.375
$ws.Range("K112").Value = 2450.0001
$ws.Range("L112").Value = 9478.125
$ws.Range("M112").Value = -1342.0001
$ws.Range("N112").Value = -11694.125
$ws.Range("H127").Value = 1091.6154
$ws.Range("I127").Value = 615.1667
$ws.Range("K127").Value = 1845.5001
$ws.Range("M127").Value = 3114.4999
$ws.Range("H138").Value = 2023
$ws.Range("I138").Value = 1539.7
$ws.Range("J138").Value = 3231.25
$ws.Range("K138").Value = 4619.1
$ws.Range("L138").Value = 9693.75
$ws.Range("M138").Value = 520.8999999999996
$ws.Range("N138").Value = -19973.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3656.9167
$ws.Range("I2").Value = 2292.875
$ws.Range("J2").Value = 6385
$ws.Range("K2").Value = 2292.875
$ws.Range("L2").Value = 6385
$ws.Range("M2").Value = -2179.875
$ws.Range("N2").Value = -6611
$ws.Range("H32").Value = 6840.65
$ws.Range("I32").Value = 7336.4443
$ws.Range("J32").Value = 2378.5
$ws.Range("K32").Value = 7336.4443
$ws.Range("L32").Value = 2378.5
$ws.Range("M32").Value = -7049.4443
$ws.Range("N32").Value = -2952.5
$ws.Range("H61").Value = 3314.9614
$ws.Range("I61").Value = 2800.625
$ws.Range("K61").Value = 2800.625
$ws.Range("M61").Value = -2588.625
$ws.Range("H74").Value = 2568.25
$ws.Range("I74").Value = 2929.2
$ws.Range("J74").Value = 1966.6666
$ws.Range("K74").Value = 2929.2
$ws.Range("L74").Value = 1966.6666
$ws.Range("M74").Value = -2055.2
$ws.Range("N74").Value = -3714.6666
$ws.Range("H77").Value = 2568.25
$ws.Range("I77").Value = 2929.2
$ws.Range("J77").Value = 1966.6666
$ws.Range("K77").Value = 14646
$ws.Range("L77").Value = 9833.333000000001
$ws.Range("M77").Value = -10278
$ws.Range("N77").Value = -18569.333
$ws.Range("H88").Value = 2280
$ws.Range("I88").Value = 2197.6
$ws.Range("J88").Value = 2362.4
$ws.Range("K88").Value = 2197.6
$ws.Range("L88").Value = 2362.4
$ws.Range("M88").Value = -1791.6
$ws.Range("N88").Value = -3174.4
$ws.Range("H91").Value = 2280
$ws.Range("I91").Value = 2197.6
$ws.Range("J91").Value = 2362.4
$ws.Range("K91").Value = 2197.6
$ws.Range("L91").Value = 2362.4
$ws.Range("M91").Value = -793.5999999999999
$ws.Range("N91").Value = -5170.4
$ws.Range("H110").Value = 2728.75
$ws.Range("I110").Value = 1087.2727
$ws.Range("K110").Value = 1087.2727
$ws.Range("M110").Value = 957.7273
$ws.Range("H116").Value = 3656.9167
$ws.Range("I116").Value = 2292.875
$ws.Range("J116").Value = 6385
$ws.Range("K116").Value = 2292.875
$ws.Range("L116").Value = 6385
$ws.Range("M116").Value = 1.125
$ws.Range("N116").Value = -10973
$ws.Range("H122").Value = 1598.6666
$ws.Range("I122").Value = 1189.2195
$ws.Range("J122").Value = 3996.8572
$ws.Range("K122").Value = 3567.6585
$ws.Range("L122").Value = 11990.5716
$ws.Range("M122").Value = -1117.6585
$ws.Range("N122").Value = -16890.5716
$ws.Range("H136").Value = 3314.9614
$ws.Range("I136").Value = 2800.625
$ws.Range("K136").Value = 8401.875
$ws.Range("M136").Value = -5851.875
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3656.9167
$ws.Range("I3").Value = 2292.875
$ws.Range("J3").Value = 6385
$ws.Range("K3").Value = 2292.875
$ws.Range("L3").Value = 6385
$ws.Range("M3").Value = -2178.875
$ws.Range("N3").Value = -6613
$ws.Range("H35").Value = 44495.5
$ws.Range("J35").Value = 44495.5
$ws.Range("L35").Value = 44495.5
$ws.Range("N35").Value = -45115.5
$ws.Range("H86").Value = 1724.5834
$ws.Range("I86").Value = 1699.5454
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1699.5454
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -576.5454
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1724.5834
$ws.Range("I89").Value = 1699.5454
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 8497.726999999999
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -2881.726999999999
$ws.Range("N89").Value = -21232
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080
$ws.Range("H134").Value = 5958.3057
$ws.Range("I134").Value = 3109.55
$ws.Range("J134").Value = 9519.25
$ws.Range("K134").Value = 9328.650000000001
$ws.Range("L134").Value = 28557.75
$ws.Range("M134").Value = -6793.650000000001
$ws.Range("N134").Value = -33627.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5557725.5
$ws.Range("I31").Value = 1488.7778
$ws.Range("J31").Value = 13892081
$ws.Range("K31").Value = 1488.7778
$ws.Range("L31").Value = 13892081
$ws.Range("M31").Value = -1193.7778
$ws.Range("N31").Value = -13892671
$ws.Range("H34").Value = 5557725.5
$ws.Range("I34").Value = 1488.7778
$ws.Range("J34").Value = 13892081
$ws.Range("K34").Value = 1488.7778
$ws.Range("L34").Value = 13892081
$ws.Range("M34").Value = -1286.7778
$ws.Range("N34").Value = -13892485
$ws.Range("H58").Value = 1728890.9
$ws.Range("I58").Value = 3825.6155
$ws.Range("J58").Value = 3130506.2
$ws.Range("K58").Value = 3825.6155
$ws.Range("L58").Value = 3130506.2
$ws.Range("M58").Value = -3622.6155
$ws.Range("N58").Value = -3130912.2
$ws.Range("H108").Value = 19796.666
$ws.Range("J108").Value = 19796.666
$ws.Range("L108").Value = 19796.666
$ws.Range("N108").Value = -27476.666
$ws.Range("H132").Value = 2899.7144
$ws.Range("I132").Value = 1576.3334
$ws.Range("J132").Value = 4664.222
$ws.Range("K132").Value = 4729.0002
$ws.Range("L132").Value = 13992.666
$ws.Range("M132").Value = -2199.0002
$ws.Range("N132").Value = -19052.666
$ws.Range("H136").Value = 1728890.9
$ws.Range("I136").Value = 3825.6155
$ws.Range("J136").Value = 3130506.2
$ws.Range("K136").Value = 11476.8465
$ws.Range("L136").Value = 9391518.600000001
$ws.Range("M136").Value = -8926.8465
$ws.Range("N136").Value = -9396618.600000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3754.4614
$ws.Range("I137").Value = 2338.2354
$ws.Range("J137").Value = 6429.5557
$ws.Range("K137").Value = 7014.706200000001
$ws.Range("L137").Value = 19288.6671
$ws.Range("M137").Value = -1914.706200000001
$ws.Range("N137").Value = -29488.6671
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3681
$ws.Range("I132").Value = 3934
$ws.Range("J132").Value = 3428
$ws.Range("K132").Value = 11802
$ws.Range("L132").Value = 10284
$ws.Range("M132").Value = -9272
$ws.Range("N132").Value = -15344
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 62503480
$ws.Range("I132").Value = 90912240
$ws.Range("J132").Value = 4199.4
$ws.Range("K132").Value = 272736720
$ws.Range("L132").Value = 12598.2
$ws.Range("M132").Value = -272734190
$ws.Range("N132").Value = -17658.2
$ws.Range("H136").Value = 33334766
$ws.Range("I136").Value = 41667924
$ws.Range("J136").Value = 2130
$ws.Range("K136").Value = 125003772
$ws.Range("L136").Value = 6390
$ws.Range("M136").Value = -125001222
$ws.Range("N136").Value = -11490
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2477.4473
$ws.Range("I132").Value = 2078.4285
$ws.Range("J132").Value = 2970.353
$ws.Range("K132").Value = 6235.2855
$ws.Range("L132").Value = 8911.059000000001
$ws.Range("M132").Value = -3705.2855
$ws.Range("N132").Value = -13971.059
